# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" sheets, per the source data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 88
$ws1.Range("F7").Value = 2666
$ws1.Range("F8").Value = 1166
$ws1.Range("F9").Value = 255
$ws1.Range("F10").Value = 110
$ws1.Range("F11").Value = 9850
$ws1.Range("F13").Value = 254
$ws1.Range("F15").Value = 609
$ws1.Range("F16").Value = 11723
$ws1.Range("F17").Value = 12052
$ws1.Range("F19").Value = 92
$ws1.Range("F21").Value = 24

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 88
$ws4.Range("F7").Value = 2666
$ws4.Range("F9").Value = 1166
$ws4.Range("F10").Value = 255
$ws4.Range("F11").Value = 110
$ws4.Range("F12").Value = 9850
$ws4.Range("F14").Value = 254
$ws4.Range("F16").Value = 609
$ws4.Range("F17").Value = 11723
$ws4.Range("F18").Value = 12052
$ws4.Range("F20").Value = 92
$ws4.Range("F22").Value = 24
